$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the current row 341, shifting the existing
# rows 341-392 down to 344-395 (carrying their formatting, e.g. the
# date style on column D, along with them).
$ws.Rows.Item(341).Insert()
$ws.Rows.Item(341).Insert()
$ws.Rows.Item(341).Insert()

# Populate the 3 newly inserted rows with the new weekly records.
# Columns A, B, C, E, F, G, R are constant across every data row in
# this sheet.

# Row 341
$ws.Cells.Item(341,1).Value = 6
$ws.Cells.Item(341,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(341,3).Value = "Metropolitana"
$ws.Cells.Item(341,4).Value = 44504
$ws.Cells.Item(341,5).Value = 13
$ws.Cells.Item(341,6).Value = 100112003
$ws.Cells.Item(341,7).Value = "Ajo"
$ws.Cells.Item(341,8).Value = "Chino"
$ws.Cells.Item(341,9).Value = "1a nueva(o)"
$ws.Cells.Item(341,10).Value = 12000
$ws.Cells.Item(341,11).Value = 2200
$ws.Cells.Item(341,12).Value = 2200
$ws.Cells.Item(341,13).Value = 2200
$ws.Cells.Item(341,14).Value = "$/paquete 20 unidades (volumen en unidades)"
$ws.Cells.Item(341,15).Value = "Región de O'Higgins"
$ws.Cells.Item(341,16).Value = 110
$ws.Cells.Item(341,17).Value = 20
$ws.Cells.Item(341,18).Value = "Hortaliza"

# Row 342
$ws.Cells.Item(342,1).Value = 6
$ws.Cells.Item(342,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(342,3).Value = "Metropolitana"
$ws.Cells.Item(342,4).Value = 44504
$ws.Cells.Item(342,5).Value = 13
$ws.Cells.Item(342,6).Value = 100112003
$ws.Cells.Item(342,7).Value = "Ajo"
$ws.Cells.Item(342,8).Value = "Chino"
$ws.Cells.Item(342,9).Value = "2a nueva(o)"
$ws.Cells.Item(342,10).Value = 25000
$ws.Cells.Item(342,11).Value = 1400
$ws.Cells.Item(342,12).Value = 1400
$ws.Cells.Item(342,13).Value = 1400
$ws.Cells.Item(342,14).Value = "$/paquete 20 unidades (volumen en unidades)"
$ws.Cells.Item(342,15).Value = "Región de O'Higgins"
$ws.Cells.Item(342,16).Value = 70
$ws.Cells.Item(342,17).Value = 20
$ws.Cells.Item(342,18).Value = "Hortaliza"

# Row 343
$ws.Cells.Item(343,1).Value = 6
$ws.Cells.Item(343,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(343,3).Value = "Metropolitana"
$ws.Cells.Item(343,4).Value = 44504
$ws.Cells.Item(343,5).Value = 13
$ws.Cells.Item(343,6).Value = 100112003
$ws.Cells.Item(343,7).Value = "Ajo"
$ws.Cells.Item(343,8).Value = "Chino"
$ws.Cells.Item(343,9).Value = "Primera"
$ws.Cells.Item(343,10).Value = 2400
$ws.Cells.Item(343,11).Value = 16500
$ws.Cells.Item(343,12).Value = 17000
$ws.Cells.Item(343,13).Value = 16729
$ws.Cells.Item(343,14).Value = "$/caja 10 kilos"
$ws.Cells.Item(343,15).Value = "China"
$ws.Cells.Item(343,16).Value = 1673
$ws.Cells.Item(343,17).Value = 10
$ws.Cells.Item(343,18).Value = "Hortaliza"
